$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" property value (row 8, column B).
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"

# Insert a new row before row 11 ("Description") to add the "Jurisdiction"
# property (with an empty value), shifting the remaining rows down by one.
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (row 11 above got the
# default/unformatted style from Insert(), so copy formats from the row
# that is now below it).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
